$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.556.22'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '3.328.44'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.62%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '3.325.89'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.178'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.576'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.37'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000269'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '673.52'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.46%  '
$ws.Range('D15').Value = '3.876.56'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '67.647.08'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = '3.333.76'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.888'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.55%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '575.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').Value = '3.691.42'
$ws.Range('E36').Value = '  -6.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.63'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.11%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.132'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.11'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.46%  '
$ws.Range('E41').Value = '  -2.60%  '
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('D44').Value = '0.0₃0666'
$ws.Range('E44').Value = '  -3.03%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.333'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('E46').Value = '  -2.69%  '
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.59%  '
